$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new row (row 7) with the "540. Single Element in a Sorted Array" question.
$ws.Range("A7").Value = 540
$ws.Range("B7").Value = "CN/LC"
$ws.Range("C7").Value = "Single Element in a Sorted Array"
$ws.Range("D7").Value = "Java"
$ws.Range("E7").Value = "Medium"

# Match formatting used by similar rows (row 6 as a base, but with its own fill for Level column).
$ws.Range("A7").Style = $ws.Range("A6").Style
$ws.Range("B7").Style = $ws.Range("B6").Style
$ws.Range("C7").Style = $ws.Range("C6").Style
$ws.Range("D7").Style = $ws.Range("D6").Style

$ws.Range("E7").HorizontalAlignment = -4108
$ws.Range("E7").VerticalAlignment = -4160
$ws.Range("E7").Interior.ThemeColor = 7
$ws.Range("E7").Interior.TintAndShade = -0.249977111117893

# Add hyperlink on the question cell.
$ws.Hyperlinks.Add($ws.Range("C7"), "https://leetcode.com/problems/single-element-in-a-sorted-array/", "", "", "Single Element in a Sorted Array")

# Update the active cell selection to match what was recorded when the workbook was saved.
$ws.Range("C13").Select()

$wb.Save()
